# Update "Última actualización" timestamps and the scraped schedule data
# for línea 141 across all three worksheets.

$wb = $excel.ActiveWorkbook

$oldTime = "02:06:16"
$newTime = "02:38:37"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"

$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 24

$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 70

$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "04:01"
$ws1.Range("D8").Value = 83

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
